# Insert a new data row at row 179 (pushes existing rows 179-226 down to 180-227)
# and populate it with a new weekly price record, matching the commit
# "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(179).Insert()

$ws.Range("A179").Value2 = 7
$ws.Range("B179").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C179").Value2 = "Ñuble"
$ws.Range("D179").Value2 = 44754
$ws.Range("E179").Value2 = 16
$ws.Range("F179").Value2 = 100112032
$ws.Range("G179").Value2 = "Zapallo italiano"
$ws.Range("H179").Value2 = "Sin especificar"
$ws.Range("I179").Value2 = "Primera"
$ws.Range("J179").Value2 = 80
$ws.Range("K179").Value2 = 11500
$ws.Range("L179").Value2 = 12500
$ws.Range("M179").Value2 = 12000
$ws.Range("N179").Value2 = "$/caja 50 unidades"
$ws.Range("O179").Value2 = "Región de Arica y Parinacota"
$ws.Range("P179").Value2 = 240
$ws.Range("Q179").Value2 = 50
$ws.Range("R179").Value2 = "Hortaliza"
